# Insert a new data row at row 86 (pushing existing rows 86-108 down to 87-109)
# and populate it with the new weekly price record for Albahaca
# (Feria Lagunitas de Puerto Montt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44642
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = 100112052
$ws.Range("G86").Value = "Albahaca"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 90
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = 5500
$ws.Range("N86").Value = '$/docena de matas'
$ws.Range("O86").Value = "Región Metropolitana"
$ws.Range("P86").Value = 917
$ws.Range("Q86").Value = 6
$ws.Range("R86").Value = "Hortaliza"
